# Add the 7th submission row (7 xgboost ensemble, 3in1 data set, 2 valid sets)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "7_281115_1041_7_xgboost_with_3in1_valid1_valid2_"
$ws.Range("B8").Value = 0.614
$ws.Range("C8").Value = "ensembled 7 xgboost, in 3in1 data set, with 2 valid sets"

# Move/restore the active cell selection as seen in the saved workbook
$ws.Range("B11").Select() | Out-Null
